# Add three new "unique" attributes (xstring_unique, xint_unique, xxref_unique)
# to the TypeTest entity: new metadata rows in the "attributes" sheet, a new
# "unique" column header, and the corresponding new data columns (with sample
# data) in the "TypeTest" sheet.

$wb = $excel.ActiveWorkbook

$wsType = $wb.Worksheets.Item("TypeTest")
$wsAttr = $wb.Worksheets.Item("attributes")

# --- order matters for how new shared strings get interned, mirroring the
# --- sequence of edits originally made in Excel ---

# 1. Header for the new unique string attribute column in TypeTest (AL1).
$wsType.Cells.Item(1, 38).Value = "xstring_unique"

# 2. New "unique" header column in the attributes metadata sheet (O1).
$wsAttr.Cells.Item(1, 15).Value = "unique"

# 3. Header for the new unique int attribute column in TypeTest (AM1).
$wsType.Cells.Item(1, 39).Value = "xint_unique"

# 4. Fill the xstring_unique column (AL) with unique string values str1..str38.
$strValues = @("str1","str2","str3","str4","str5","str6","str7","str8","str9","str10", `
    "str11","str12","str13","str14","str15","str16","str17","str18","str19","str20", `
    "str21","str22","str23","str24","str25","str26","str27","str28","str29","str30", `
    "str31","str32","str33","str34","str35","str36","str37","str38")

for ($i = 0; $i -lt $strValues.Length; $i++) {
    $wsType.Cells.Item($i + 2, 38).Value = $strValues[$i]
}

# 5. Fill the xint_unique column (AM) with unique integer values 1..38.
for ($i = 1; $i -le 38; $i++) {
    $wsType.Cells.Item($i + 1, 39).Value = $i
}

# 6. Header for the new unique xref attribute column in TypeTest (AN1), and
#    fill the first five rows with unique references (further rows would
#    repeat references and violate uniqueness, so they stay blank).
$wsType.Cells.Item(1, 40).Value = "xxref_unique"

$xrefValues = @("ref1","ref2","ref3","ref4","ref5")
for ($i = 0; $i -lt $xrefValues.Length; $i++) {
    $wsType.Cells.Item($i + 2, 40).Value = $xrefValues[$i]
}

# --- new attribute metadata rows describing the three new columns ---

# xstring_unique
$wsAttr.Cells.Item(41, 1).Value = "xstring_unique"
$wsAttr.Cells.Item(41, 2).Value = "TypeTest"
$wsAttr.Cells.Item(41, 3).Value = "string"
$wsAttr.Cells.Item(41, 5).Value = $false
$wsAttr.Cells.Item(41, 6).Value = $false
$wsAttr.Cells.Item(41, 15).Value = $true

# xint_unique
$wsAttr.Cells.Item(42, 1).Value = "xint_unique"
$wsAttr.Cells.Item(42, 2).Value = "TypeTest"
$wsAttr.Cells.Item(42, 3).Value = "int"
$wsAttr.Cells.Item(42, 5).Value = $false
$wsAttr.Cells.Item(42, 6).Value = $false
$wsAttr.Cells.Item(42, 15).Value = $true

# xxref_unique
$wsAttr.Cells.Item(43, 1).Value = "xxref_unique"
$wsAttr.Cells.Item(43, 2).Value = "TypeTest"
$wsAttr.Cells.Item(43, 3).Value = "xref"
$wsAttr.Cells.Item(43, 4).Value = "TypeTestRef"
$wsAttr.Cells.Item(43, 5).Value = $false
$wsAttr.Cells.Item(43, 6).Value = $true
$wsAttr.Cells.Item(43, 15).Value = $true
